$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '28.072.27'
$ws.Range("E2").Value = '  -0.38%  '
$ws.Range("D3").Value = '1.876.83'
$ws.Range("E3").Value = '  -1.64%  '
$ws.Range("E4").Value = '  +0.35%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '313.58'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.19%  '
$ws.Range("E6").Value = '  +0.27%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.5088'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +0.10%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3855'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -2.04%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.09164'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -1.76%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '1.124'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -1.54%  '
$ws.Range("B11").Value = 'Polkadot'
$ws.Range("C11").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '6.353'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -0.65%  '
$ws.Range("B12").Value = 'Solana'
$ws.Range("C12").Value = 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '20.74'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -0.88%  '
$ws.Range("B13").Value = 'WrappedEther'
$ws.Range("C13").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D13").Value = '1.876.72'
$ws.Range("E13").Value = '  -1.19%  '
$ws.Range("B14").Value = 'Chainlink'
$ws.Range("C14").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '7.202'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -1.51%  '
$ws.Range("B15").Value = 'BinanceUSD'
$ws.Range("C15").Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '1.002'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +0.39%  '
$ws.Range("B16").Value = 'ShibaInu'
$ws.Range("C16").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.00001117'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -0.62%  '
$ws.Range("B17").Value = 'Litecoin'
$ws.Range("C17").Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '91.11'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -1.55%  '
$ws.Range("B18").Value = 'TRON'
$ws.Range("C18").Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.06598'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +0.20%  '
$ws.Range("B19").Value = 'Avalanche'
$ws.Range("C19").Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '18.19'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +1.47%  '
$ws.Range("B20").Value = 'Dai'
$ws.Range("C20").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '1.002'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +0.37%  '
$ws.Range("B21").Value = 'Uniswap'
$ws.Range("C21").Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '6.119'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -1.76%  '
$ws.Range("B22").Value = 'WrappedBTC'
$ws.Range("C22").Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range("D22").Value = '28.091.06'
$ws.Range("E22").Value = '  -0.45%  '
$ws.Range("B23").Value = 'Cosmos'
$ws.Range("C23").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '11.43'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +0.33%  '
$ws.Range("B24").Value = 'Toncoin'
$ws.Range("C24").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.278'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -1.65%  '
$ws.Range("B25").Value = 'LidoDAOToken'
$ws.Range("C25").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.546'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -3.68%  '
$ws.Range("B26").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C26").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D26").Value = '2.089.72'
$ws.Range("E26").Value = '  -1.66%  '
$ws.Range("B27").Value = 'EthereumClassic'
$ws.Range("C27").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '20.82'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -1.15%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '157.71'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +0.13%  '
$ws.Range("B29").Value = 'BitcoinCash'
$ws.Range("C29").Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '126.69'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -0.50%  '
$ws.Range("B30").Value = 'Stellar'
$ws.Range("C30").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.1053'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -1.67%  '
$ws.Range("B31").Value = 'ImmutableX'
$ws.Range("C31").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.064'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -2.73%  '
$ws.Range("B32").Value = 'Filecoin'
$ws.Range("C32").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '5.617'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -0.17%  '
$ws.Range("B33").Value = 'HuobiToken'
$ws.Range("C33").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '3.603'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -0.11%  '
$ws.Range("B34").Value = 'FraxShare'
$ws.Range("C34").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '9.673'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +0.20%  '
$ws.Range("B35").Value = 'Hedera'
$ws.Range("C35").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.06578'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -1.65%  '
$ws.Range("B36").Value = 'VeChain'
$ws.Range("C36").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.02429'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +0.64%  '
$ws.Range("B37").Value = 'Algorand'
$ws.Range("C37").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.2181'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -0.26%  '
$ws.Range("B38").Value = 'ARBITRUM'
$ws.Range("C38").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.211'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -2.35%  '
$ws.Range("B39").Value = 'TrustWalletToken'
$ws.Range("C39").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '1.267'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +0.30%  '
$ws.Range("B40").Value = 'TheSandbox'
$ws.Range("C40").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.6401'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -0.05%  '
$ws.Range("B41").Value = 'Aptos'
$ws.Range("C41").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '11.46'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -0.54%  '
$ws.Range("B42").Value = 'InternetComputer(DFINITY)'
$ws.Range("C42").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '4.922'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -1.54%  '
$ws.Range("B43").Value = 'Frax'
$ws.Range("C43").Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.001'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +0.29%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.6023'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +0.17%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '13.19'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -0.44%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '3.671'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -0.93%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.275'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -0.14%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.238'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +4.76%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.998'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -1.20%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '121.51'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -1.27%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '79.79'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +1.38%  '
